$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 351569
$ws.Range("D2").Value = 446602088
$ws.Range("C3").Value = 286
$ws.Range("D3").Value = 344716
$ws.Range("C10").Value = 125701
$ws.Range("D10").Value = 184076137
$ws.Range("C11").Value = 149
$ws.Range("D11").Value = 220513
$ws.Range("C12").Value = 65939
$ws.Range("D12").Value = 95174509
$ws.Range("C16").Value = 4160
$ws.Range("D16").Value = 5907077
$ws.Range("C21").Value = 8682
$ws.Range("D21").Value = 12116443
$ws.Range("C23").Value = 83685
$ws.Range("D23").Value = 103863029
$ws.Range("C29").Value = 34406
$ws.Range("D29").Value = 50348864
$ws.Range("C32").Value = 12507
$ws.Range("D32").Value = 17993728
$ws.Range("C35").Value = 1620
$ws.Range("D35").Value = 2274902
$ws.Range("C37").Value = 2293
$ws.Range("D37").Value = 3229933
$ws.Range("C38").Value = 104591
$ws.Range("D38").Value = 130963114
$ws.Range("C44").Value = 924
$ws.Range("D44").Value = 1360685
$ws.Range("C46").Value = 46894
$ws.Range("D46").Value = 68684716
$ws.Range("C47").Value = 32
$ws.Range("D47").Value = 46474
$ws.Range("C48").Value = 9976
$ws.Range("D48").Value = 14297970
$ws.Range("C50").Value = 1485
$ws.Range("D50").Value = 2067757
$ws.Range("C53").Value = 2904
$ws.Range("D53").Value = 4071147
$ws.Range("C54").Value = 74713
$ws.Range("D54").Value = 93432428
$ws.Range("C61").Value = 30157
$ws.Range("D61").Value = 44189860
$ws.Range("C64").Value = 12200
$ws.Range("D64").Value = 17624615
$ws.Range("C66").Value = 1430
$ws.Range("D66").Value = 2000166
$ws.Range("C70").Value = 1868
$ws.Range("D70").Value = 2624327
$ws.Range("C72").Value = 22494
$ws.Range("D72").Value = 29378508
$ws.Range("C76").Value = 8284
$ws.Range("D76").Value = 12130524
$ws.Range("C78").Value = 5683
$ws.Range("D78").Value = 8263931
$ws.Range("C79").Value = 556
$ws.Range("D79").Value = 784968
$ws.Range("C80").Value = 361
$ws.Range("D80").Value = 509448
$ws.Range("C81").Value = 152802
$ws.Range("D81").Value = 189674839
$ws.Range("C87").Value = 67577
$ws.Range("D87").Value = 98985778
$ws.Range("C90").Value = 32159
$ws.Range("D90").Value = 46522916
$ws.Range("C92").Value = 2985
$ws.Range("D92").Value = 4306182
$ws.Range("C94").Value = 3541
$ws.Range("D94").Value = 4998301
$ws.Range("C95").Value = 37839
$ws.Range("D95").Value = 51324684
$ws.Range("C99").Value = 9383
$ws.Range("D99").Value = 13784762
$ws.Range("C101").Value = 8725
$ws.Range("D101").Value = 12659926
$ws.Range("C103").Value = 598
$ws.Range("D103").Value = 845430
$ws.Range("C104").Value = 571
$ws.Range("D104").Value = 826560
$ws.Range("C105").Value = 16517
$ws.Range("D105").Value = 30882945
$ws.Range("C108").Value = 3831
$ws.Range("D108").Value = 7680181
$ws.Range("C110").Value = 5437
$ws.Range("D110").Value = 11023520
$ws.Range("C112").Value = 252
$ws.Range("D112").Value = 507340
$ws.Range("C113").Value = 322
$ws.Range("D113").Value = 626650
$ws.Range("C115").Value = 153790
$ws.Range("D115").Value = 189806410
$ws.Range("C121").Value = 56558
$ws.Range("D121").Value = 82820395
$ws.Range("C123").Value = 30412
$ws.Range("D123").Value = 44053367
$ws.Range("C127").Value = 2827
$ws.Range("D127").Value = 3983097
$ws.Range("C129").Value = 614664
$ws.Range("D129").Value = 811361164
$ws.Range("C131").Value = 249
$ws.Range("D131").Value = 366713
$ws.Range("C134").Value = 1593
$ws.Range("D134").Value = 2359103
$ws.Range("C136").Value = 236809
$ws.Range("D136").Value = 347877220
$ws.Range("C137").Value = 539
$ws.Range("D137").Value = 803147
$ws.Range("C139").Value = 220897
$ws.Range("D139").Value = 321208427
$ws.Range("C142").Value = 3019
$ws.Range("D142").Value = 4243734
$ws.Range("C145").Value = 8502
$ws.Range("D145").Value = 11981868
$ws.Range("C148").Value = 49534
$ws.Range("D148").Value = 66040199
$ws.Range("C149").Value = 34
$ws.Range("D149").Value = 45310
$ws.Range("C154").Value = 15278
$ws.Range("D154").Value = 22390110
$ws.Range("C155").Value = 4149
$ws.Range("D155").Value = 5987826
$ws.Range("C160").Value = 498
$ws.Range("D160").Value = 704259
$ws.Range("C161").Value = 19793
$ws.Range("D161").Value = 26168853
$ws.Range("C164").Value = 57
$ws.Range("D164").Value = 83906
$ws.Range("C165").Value = 8193
$ws.Range("D165").Value = 11919200
$ws.Range("C167").Value = 5815
$ws.Range("D167").Value = 8369992
$ws.Range("C172").Value = 28736
$ws.Range("D172").Value = 57862448
$ws.Range("C173").Value = 2867
$ws.Range("D173").Value = 5702117
$ws.Range("C178").Value = 94607
$ws.Range("D178").Value = 117784627
$ws.Range("C185").Value = 36156
$ws.Range("D185").Value = 52988636
$ws.Range("C187").Value = 14384
$ws.Range("D187").Value = 20778331
$ws.Range("C191").Value = 2038
$ws.Range("D191").Value = 2863791
$ws.Range("C193").Value = 257955
$ws.Range("D193").Value = 319488138
$ws.Range("C199").Value = 940
$ws.Range("D199").Value = 1381781
$ws.Range("C201").Value = 92351
$ws.Range("D201").Value = 135301215
$ws.Range("C204").Value = 36264
$ws.Range("D204").Value = 52210076
$ws.Range("C207").Value = 5470
$ws.Range("D207").Value = 7790817
$ws.Range("C210").Value = 6057
$ws.Range("D210").Value = 8384838
$ws.Range("C213").Value = 286946
$ws.Range("D213").Value = 354017372
$ws.Range("C220").Value = 656
$ws.Range("D220").Value = 954893
$ws.Range("C222").Value = 102132
$ws.Range("D222").Value = 149356204
$ws.Range("C225").Value = 56660
$ws.Range("D225").Value = 81851102
$ws.Range("C228").Value = 4933
$ws.Range("D228").Value = 6918902
$ws.Range("C231").Value = 7524
$ws.Range("D231").Value = 10414314
$ws.Range("C234").Value = 115714
$ws.Range("D234").Value = 144140155
$ws.Range("C241").Value = 52683
$ws.Range("D241").Value = 77149658
$ws.Range("C243").Value = 13893
$ws.Range("D243").Value = 19985996
$ws.Range("C247").Value = 3089
$ws.Range("D247").Value = 4321230
$ws.Range("C248").Value = 283587
$ws.Range("D248").Value = 357505627
$ws.Range("C249").Value = 196
$ws.Range("D249").Value = 244328
$ws.Range("C255").Value = 915
$ws.Range("D255").Value = 1341846
$ws.Range("C257").Value = 104461
$ws.Range("D257").Value = 153011220
$ws.Range("C260").Value = 73383
$ws.Range("D260").Value = 106347428
$ws.Range("C261").Value = 19
$ws.Range("D261").Value = 28489
$ws.Range("C262").Value = 2554
$ws.Range("D262").Value = 3595912
$ws.Range("C265").Value = 6000
$ws.Range("D265").Value = 8413116
